$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.9714294676670578
$ws.Range("D2").Value = 0.3418946047018387

$ws.Range("C3").Value = 0.1653503009920105
$ws.Range("D3").Value = 0.8701782837629068

$ws.Range("C4").Value = 1.109888893704867
$ws.Range("D4").Value = 0.2790316597188462

$ws.Range("C5").Value = 0.3876083195771526
$ws.Range("D5").Value = 0.7020306402195642

$ws.Range("C6").Value = 0.7544342953799321
$ws.Range("D6").Value = 0.4585885903545193

$ws.Range("C7").Value = 2.367574704637302
$ws.Range("D7").Value = 0.02711771962191434

$ws.Range("C8").Value = 1.165845893028361
$ws.Range("D8").Value = 0.2561610537089749

$ws.Range("C9").Value = 0.706255693091176
$ws.Range("D9").Value = 0.4874419189203691

$ws.Range("C10").Value = 0.2815029591027292
$ws.Range("D10").Value = 0.7809554919050159

$ws.Range("C11").Value = -0.6684658283334465
$ws.Range("D11").Value = 0.5107897776554307
